# Lattice multiplication exercises: regenerate the 15 table-cell
# expressions/partial-products in place (structure/formatting unchanged).
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1, Col 1: '32 x 58' -> '37 x 40'
$cell = $t.Cell(1, 1)
$cell.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>37 x 40</w:t><w:br/><w:t xml:space="preserve">  4    0</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>3|    |</w:t><w:br/><w:t>7|    |</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>') | Out-Null

# Row 1, Col 2: '77 x 78' -> '32 x 18'
$cell = $t.Cell(1, 2)
$cell.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>32 x 18</w:t><w:br/><w:t xml:space="preserve">  1    8</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>3|    |</w:t><w:br/><w:t>2|    |</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>') | Out-Null

# Row 1, Col 3: '26 x 12' -> '45 x 20'
$cell = $t.Cell(1, 3)
$cell.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>45 x 20</w:t><w:br/><w:t xml:space="preserve">  2    0</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>4|    |</w:t><w:br/><w:t>5|    |</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>') | Out-Null

# Row 2, Col 1: '32 x 42' -> '28 x 65'
$cell = $t.Cell(2, 1)
$cell.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>28 x 65</w:t><w:br/><w:t xml:space="preserve">  6    5</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>2|    |</w:t><w:br/><w:t>8|    |</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>') | Out-Null

# Row 2, Col 2: '97 x 67' -> '77 x 32'
$cell = $t.Cell(2, 2)
$cell.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>77 x 32</w:t><w:br/><w:t xml:space="preserve">  3    2</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>7|    |</w:t><w:br/><w:t>7|    |</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>') | Out-Null

# Row 2, Col 3: '54 x 58' -> '95 x 20'
$cell = $t.Cell(2, 3)
$cell.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>95 x 20</w:t><w:br/><w:t xml:space="preserve">  2    0</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>9|    |</w:t><w:br/><w:t>5|    |</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>') | Out-Null

# Row 3, Col 1: '61 x 86' -> '35 x 87'
$cell = $t.Cell(3, 1)
$cell.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>35 x 87</w:t><w:br/><w:t xml:space="preserve">  8    7</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>3|    |</w:t><w:br/><w:t>5|    |</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>') | Out-Null

# Row 3, Col 2: '71 x 49' -> '34 x 17'
$cell = $t.Cell(3, 2)
$cell.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>34 x 17</w:t><w:br/><w:t xml:space="preserve">  1    7</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>3|    |</w:t><w:br/><w:t>4|    |</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>') | Out-Null

# Row 3, Col 3: '98 x 67' -> '16 x 44'
$cell = $t.Cell(3, 3)
$cell.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>16 x 44</w:t><w:br/><w:t xml:space="preserve">  4    4</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>1|    |</w:t><w:br/><w:t>6|    |</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>') | Out-Null

# Row 4, Col 1: '79 x 15' -> '10 x 94'
$cell = $t.Cell(4, 1)
$cell.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>10 x 94</w:t><w:br/><w:t xml:space="preserve">  9    4</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>1|    |</w:t><w:br/><w:t>0|    |</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>') | Out-Null

# Row 4, Col 2: '56 x 20' -> '10 x 72'
$cell = $t.Cell(4, 2)
$cell.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>10 x 72</w:t><w:br/><w:t xml:space="preserve">  7    2</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>1|    |</w:t><w:br/><w:t>0|    |</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>') | Out-Null

# Row 4, Col 3: '64 x 72' -> '88 x 29'
$cell = $t.Cell(4, 3)
$cell.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>88 x 29</w:t><w:br/><w:t xml:space="preserve">  2    9</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>8|    |</w:t><w:br/><w:t>8|    |</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>') | Out-Null

# Row 5, Col 1: '66 x 89' -> '68 x 48'
$cell = $t.Cell(5, 1)
$cell.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>68 x 48</w:t><w:br/><w:t xml:space="preserve">  4    8</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>6|    |</w:t><w:br/><w:t>8|    |</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>') | Out-Null

# Row 5, Col 2: '87 x 52' -> '19 x 50'
$cell = $t.Cell(5, 2)
$cell.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>19 x 50</w:t><w:br/><w:t xml:space="preserve">  5    0</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>1|    |</w:t><w:br/><w:t>9|    |</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>') | Out-Null

# Row 5, Col 3: '91 x 78' -> '63 x 95'
$cell = $t.Cell(5, 3)
$cell.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>63 x 95</w:t><w:br/><w:t xml:space="preserve">  9    5</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>6|    |</w:t><w:br/><w:t>3|    |</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>') | Out-Null

Write-Output "Updated lattice multiplication exercises (15 cells)."